$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new row above row 14 ("notify_prelaunch ..." etc. shift
# down by one row, 11:19 commands staying intact, new blank row 14
# is created for the new "reboot" command).
# ------------------------------------------------------------------
$ws.Rows("14:14").Insert()

# The blank row created by Insert doesn't fully inherit every cell's
# formatting, so pull full formatting from row 13 (identical template
# row) into the new row 14.
$ws.Range("A13:AC13").Copy()
$ws.Range("A14:AC14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Populate new row 14 with the "reboot" timepix command.
# ------------------------------------------------------------------
$ws.Range("A14").Value = "reboot"
$ws.Range("B14").Value = "0000 0000"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 1111001
$ws.Range("F14").Value = [char]0x2014
$ws.Range("G14").Value = [char]0x2014
$ws.Range("H14").Value = [char]0x2014
$ws.Range("I14").Value = [char]0x2014
$ws.Range("J14").Value = [char]0x2014
$ws.Range("K14").Value = "0x00"
$ws.Range("L14:V14").Value = 0
$ws.Range("W14").Value = 1
$ws.Range("X14").Value = "?"
$ws.Range("Y14").Value = "0x00"
$ws.Range("Z14").Value = "0x00"
$ws.Range("AB14").Value = "todo"
$ws.Range("AC14").Value = "reboot raspberry pi"

# ------------------------------------------------------------------
# Re-apply the hex/address shared formulas across the full command
# range so that rows 11-18 (and the new 19-20 duplicate block) keep
# computing correctly after the insert.
# ------------------------------------------------------------------
$ws.Range("E11:E18").Formula = "=_xlfn.CONCAT(""0x"", DEC2HEX(_xlfn.BITLSHIFT(`$C11,7) + BIN2DEC(`$D11)))"
$ws.Range("E19:E20").Formula = "=_xlfn.CONCAT(""0x"", DEC2HEX(_xlfn.BITLSHIFT(`$C19,7) + BIN2DEC(`$D19)))"
$ws.Range("AA11:AA20").Formula = "=`$E11"

# ------------------------------------------------------------------
# Update sheet view: frozen-pane scroll position & active selection.
# ------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 17
$win.ScrollRow = 3
$ws.Range("A14").Select()

Write-Host "edit complete"
